# Fix header row 2: replace placeholder "unnamedXXX" labels with "total"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
